$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2022 / 2023 year columns to the header row
$ws.Range("T1").Value = 2022
$ws.Range("U1").Value = 2023

# Add the new "additional_accounts" (covid-19) values for 2022 and 2023
$ws.Range("T9").Value = 4284
$ws.Range("U9").Value = 1179

# Match the thousands-separator number formatting used by the other
# "additional_accounts" cells in that row (e.g. S9)
$ws.Range("T9:U9").NumberFormat = "#,##0"

# Move the active selection to reflect where editing continued
$ws.Range("N5").Select() | Out-Null
